# Updated PDF Merger to latest version
# Rename the aspect-score headers to aspect-intensity headers, and add a
# running "Aspect intensity total" column (M) that sums the previous row's
# Sun/Moon/Asc aspect columns (J:L) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Katie Scruggs")

# Rename headers in row 1
$ws.Range("J1").Value = "Sun Aspect Intensity"
$ws.Range("K1").Value = "Moon Aspect Intensity"
$ws.Range("L1").Value = "Asc Aspect Intensity"
$ws.Range("M1").Value = "Aspect intensity total"

# Add the new M column formulas for every data row (2 through 121), each
# one summing the J:L cells of the row directly above it.
$lastRow = 121
for ($row = 2; $row -le $lastRow; $row++) {
    $prevRow = $row - 1
    $ws.Range("M$row").Formula = "=SUM(J$prevRow`:L$prevRow)"
}
